$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.660.21'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '2.293.38'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '114.77'
$ws.Range("E5").Value = '  +19.35%  '
$ws.Range("D6").Value = '268.67'
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").Value = '0.622'
$ws.Range("E9").Value = '  +2.03%  '
$ws.Range("D10").Value = '48.25'
$ws.Range("E10").Value = '  +5.29%  '
$ws.Range("D12").Value = '8.79'
$ws.Range("E12").Value = '  +12.72%  '
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").Value = '15.62'
$ws.Range("E14").Value = '  +3.20%  '
$ws.Range("D15").Value = '2.634.11'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '0.848'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '2.286.78'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '43.671.83'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("E19").Value = '  +2.42%  '
$ws.Range("D20").Value = '6.52'
$ws.Range("E20").Value = '  +5.07%  '
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("E22").Value = '  -1.44%  '
$ws.Range("D23").Value = '233.21'
$ws.Range("D24").Value = '9.76'
$ws.Range("E24").Value = '  +6.44%  '
$ws.Range("D25").Value = '2.82'
$ws.Range("E25").Value = '  +12.73%  '
$ws.Range("D27").Value = '11.61'
$ws.Range("E27").Value = '  +4.26%  '
$ws.Range("D28").Value = '42.06'
$ws.Range("E28").Value = '  +3.78%  '
$ws.Range("E29").Value = '  -2.03%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = '176.61'
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").Value = '0.0935'
$ws.Range("E32").Value = '  +4.94%  '
$ws.Range("D33").Value = '21.60'
$ws.Range("E33").Value = '  -0.96%  '
$ws.Range("D34").Value = '5.58'
$ws.Range("E34").Value = '  +4.17%  '
$ws.Range("E35").Value = '  +1.05%  '
$ws.Range("D36").Value = '4.74'
$ws.Range("E36").Value = '  +9.42%  '
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").Value = '0.0357'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("E39").Value = '  +12.57%  '
$ws.Range("D40").Value = '2.44'
$ws.Range("E40").Value = '  +5.63%  '
$ws.Range("D41").Value = '13.87'
$ws.Range("E41").Value = '  +13.23%  '
$ws.Range("E42").Value = '  +3.52%  '
$ws.Range("D43").Value = '72.94'
$ws.Range("E43").Value = '  +11.53%  '
$ws.Range("E44").Value = '  +7.28%  '
$ws.Range("D45").Value = '6.11'
$ws.Range("E45").Value = '  +17.16%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '8.73'
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("D48").Value = '102.86'
$ws.Range("E48").Value = '  +5.87%  '
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("E50").Value = '  +3.49%  '
$ws.Range("E51").Value = '  +3.84%  '
